$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of price data was inserted at the top of the historical block
# (rows 214-257), pushing the existing rows down by two (to 216-259).
$ws.Rows("214:215").Insert()

# Row 214 - "Primera" quality entry for the new week
$ws.Range("A214").Value = 8
$ws.Range("B214").Value = "Terminal La Palmera de La Serena"
$ws.Range("C214").Value = "Coquimbo"
$ws.Range("D214").Value = 44637
$ws.Range("E214").Value = 4
$ws.Range("F214").Value = 100114014
$ws.Range("G214").Value = "Betarraga"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 2200
$ws.Range("K214").Value = 500
$ws.Range("L214").Value = 600
$ws.Range("M214").Value = 550
$ws.Range("N214").Value = '$/paquete 3 unidades'
$ws.Range("O214").Value = "Provincia del Elquí"
$ws.Range("P214").Value = 183
$ws.Range("Q214").Value = 3
$ws.Range("R214").Value = "Hortaliza"

# Row 215 - "Segunda" quality entry for the new week
$ws.Range("A215").Value = 8
$ws.Range("B215").Value = "Terminal La Palmera de La Serena"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 44637
$ws.Range("E215").Value = 4
$ws.Range("F215").Value = 100114014
$ws.Range("G215").Value = "Betarraga"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Segunda"
$ws.Range("J215").Value = 1440
$ws.Range("K215").Value = 400
$ws.Range("L215").Value = 450
$ws.Range("M215").Value = 425
$ws.Range("N215").Value = '$/paquete 3 unidades'
$ws.Range("O215").Value = "Provincia del Elquí"
$ws.Range("P215").Value = 142
$ws.Range("Q215").Value = 3
$ws.Range("R215").Value = "Hortaliza"
